$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 776.6177
$ws.Range("J17").Value = 776.6177
$ws.Range("L17").Value = 2329.8531
$ws.Range("N17").Value = -2665.8531
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
# Row 33
$ws.Range("H33").Value = 527.381
$ws.Range("I33").Value = 98.5
$ws.Range("J33").Value = 1899.8
$ws.Range("K33").Value = 98.5
$ws.Range("L33").Value = 1899.8
$ws.Range("M33").Value = 130.5
$ws.Range("N33").Value = -2357.8
# Row 69
$ws.Range("H69").Value = 2970.8333
$ws.Range("I69").Value = 2500
$ws.Range("J69").Value = 3013.6365
$ws.Range("K69").Value = 7500
$ws.Range("L69").Value = 9040.9095
$ws.Range("M69").Value = -6626
$ws.Range("N69").Value = -10788.9095
# Row 72
$ws.Range("H72").Value = 2970.8333
$ws.Range("I72").Value = 2500
$ws.Range("J72").Value = 3013.6365
$ws.Range("K72").Value = 22500
$ws.Range("L72").Value = 27122.7285
$ws.Range("M72").Value = -18132
$ws.Range("N72").Value = -35858.7285
# Row 107
$ws.Range("H107").Value = 183.95
$ws.Range("I107").Value = 183.95
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 183.95
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1736.05
$ws.Range("N107").ClearContents()
# Row 113
$ws.Range("H113").Value = 4119.067
$ws.Range("I113").Value = 3735
$ws.Range("J113").Value = 4558
$ws.Range("K113").Value = 3735
$ws.Range("L113").Value = 4558
$ws.Range("M113").Value = -481
$ws.Range("N113").Value = -11066
# Row 129
$ws.Range("H129").Value = 914.8148
$ws.Range("I129").Value = 295
$ws.Range("J129").Value = 1055.6818
$ws.Range("K129").Value = 885
$ws.Range("L129").Value = 3167.0454
$ws.Range("M129").Value = 4115
$ws.Range("N129").Value = -13167.0454
# Row 132
$ws.Range("H132").Value = 861639.25
$ws.Range("I132").Value = 1846.2693
$ws.Range("J132").Value = 9803486
$ws.Range("K132").Value = 5538.8079
$ws.Range("L132").Value = 29410458
$ws.Range("M132").Value = -3008.8079
$ws.Range("N132").Value = -29415518
# Row 135
$ws.Range("H135").Value = 23901.87
$ws.Range("I135").Value = 31014
$ws.Range("J135").Value = 3750.8333
$ws.Range("K135").Value = 279126
$ws.Range("L135").Value = 33757.4997
$ws.Range("M135").Value = -276591
$ws.Range("N135").Value = -38827.4997
# Row 137
$ws.Range("H137").Value = 2382661.8
$ws.Range("I137").Value = 3704932
$ws.Range("J137").Value = 2575
$ws.Range("K137").Value = 11114796
$ws.Range("L137").Value = 7725
$ws.Range("M137").Value = -11112246
$ws.Range("N137").Value = -12825
# Row 138
$ws.Range("H138").Value = 2452992
$ws.Range("I138").Value = 1417.6305
$ws.Range("J138").Value = 7579011.5
$ws.Range("K138").Value = 4252.8915
$ws.Range("L138").Value = 22737034.5
$ws.Range("M138").Value = 887.1085000000003
$ws.Range("N138").Value = -22747314.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1063.64
$ws.Range("I32").Value = 821.34784
$ws.Range("J32").Value = 3850
$ws.Range("K32").Value = 821.34784
$ws.Range("L32").Value = 3850
$ws.Range("M32").Value = -534.34784
$ws.Range("N32").Value = -4424
# Row 61
$ws.Range("H61").Value = 28630160
$ws.Range("I61").Value = 38501430
$ws.Range("J61").Value = 113158.555
$ws.Range("K61").Value = 38501430
$ws.Range("L61").Value = 113158.555
$ws.Range("M61").Value = -38501218
$ws.Range("N61").Value = -113582.555
# Row 74
$ws.Range("H74").Value = 7799770.5
$ws.Range("I74").Value = 9834399
$ws.Range("J74").Value = 113395.555
$ws.Range("K74").Value = 9834399
$ws.Range("L74").Value = 113395.555
$ws.Range("M74").Value = -9833525
$ws.Range("N74").Value = -115143.555
# Row 77
$ws.Range("H77").Value = 7799770.5
$ws.Range("I77").Value = 9834399
$ws.Range("J77").Value = 113395.555
$ws.Range("K77").Value = 49171995
$ws.Range("L77").Value = 566977.7749999999
$ws.Range("M77").Value = -49167627
$ws.Range("N77").Value = -575713.7749999999
# Row 97
$ws.Range("H97").Value = 1839287.1
$ws.Range("I97").Value = 2605305.8
$ws.Range("J97").Value = 842.2
$ws.Range("K97").Value = 2605305.8
$ws.Range("L97").Value = 842.2
$ws.Range("M97").Value = -2604809.8
$ws.Range("N97").Value = -1834.2
# Row 132
$ws.Range("H132").Value = 39690.85
$ws.Range("I132").Value = 30148.795
$ws.Range("J132").Value = 55912.35
$ws.Range("K132").Value = 90446.385
$ws.Range("L132").Value = 167737.05
$ws.Range("M132").Value = -87916.385
$ws.Range("N132").Value = -172797.05
# Row 136
$ws.Range("H136").Value = 28630160
$ws.Range("I136").Value = 38501430
$ws.Range("J136").Value = 113158.555
$ws.Range("K136").Value = 115504290
$ws.Range("L136").Value = 339475.665
$ws.Range("M136").Value = -115501740
$ws.Range("N136").Value = -344575.665

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 623.36365
$ws.Range("I94").Value = 405.875
$ws.Range("J94").Value = 1203.3334
$ws.Range("K94").Value = 405.875
$ws.Range("L94").Value = 1203.3334
$ws.Range("M94").Value = 45.125
$ws.Range("N94").Value = -2105.3334
# Row 99
$ws.Range("H99").Value = 1132.625
$ws.Range("I99").Value = 960
$ws.Range("J99").Value = 1420.3334
$ws.Range("K99").Value = 960
$ws.Range("L99").Value = 1420.3334
$ws.Range("M99").Value = 538
$ws.Range("N99").Value = -4416.3334
# Row 105
$ws.Range("H105").Value = 20002346
$ws.Range("I105").Value = 33335630
$ws.Range("J105").Value = 2420
$ws.Range("K105").Value = 33335630
$ws.Range("L105").Value = 2420
$ws.Range("M105").Value = -33333883
$ws.Range("N105").Value = -5914
# Row 134
$ws.Range("H134").Value = 1421.2167
$ws.Range("I134").Value = 894.75
$ws.Range("J134").Value = 3527.0833
$ws.Range("K134").Value = 2684.25
$ws.Range("L134").Value = 10581.2499
$ws.Range("M134").Value = -149.25
$ws.Range("N134").Value = -15651.2499

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 166667680
$ws.Range("I16").Value = 2011
$ws.Range("K16").Value = 2011
$ws.Range("M16").Value = -1724
# Row 31
$ws.Range("H31").Value = 3440.853
$ws.Range("I31").Value = 1386.4286
$ws.Range("K31").Value = 1386.4286
$ws.Range("M31").Value = -1091.4286
# Row 34
$ws.Range("H34").Value = 3440.853
$ws.Range("I34").Value = 1386.4286
$ws.Range("K34").Value = 1386.4286
$ws.Range("M34").Value = -1184.4286
# Row 58
$ws.Range("H58").Value = 20409842
$ws.Range("I58").Value = 25642536
$ws.Range("J58").Value = 2338.8
$ws.Range("K58").Value = 25642536
$ws.Range("L58").Value = 2338.8
$ws.Range("M58").Value = -25642333
$ws.Range("N58").Value = -2744.8
# Row 113
$ws.Range("H113").Value = 166667680
$ws.Range("I113").Value = 2011
$ws.Range("K113").Value = 2011
$ws.Range("M113").Value = 159
# Row 122
$ws.Range("H122").Value = 2487.5386
$ws.Range("I122").Value = 2003
$ws.Range("J122").Value = 3262.8
$ws.Range("K122").Value = 6009
$ws.Range("L122").Value = 9788.400000000001
$ws.Range("M122").Value = -3559
$ws.Range("N122").Value = -14688.4
# Row 132
$ws.Range("H132").Value = 21955.918
$ws.Range("I132").Value = 1359.7949
$ws.Range("J132").Value = 102280.8
$ws.Range("K132").Value = 4079.384700000001
$ws.Range("L132").Value = 306842.4
$ws.Range("M132").Value = -1549.384700000001
$ws.Range("N132").Value = -311902.4
# Row 134
$ws.Range("H134").Value = 19508.406
$ws.Range("I134").Value = 1189.8292
$ws.Range("J134").Value = 61234.055
$ws.Range("K134").Value = 3569.487599999999
$ws.Range("L134").Value = 183702.165
$ws.Range("M134").Value = -1034.487599999999
$ws.Range("N134").Value = -188772.165
# Row 136
$ws.Range("H136").Value = 20409842
$ws.Range("I136").Value = 25642536
$ws.Range("J136").Value = 2338.8
$ws.Range("K136").Value = 76927608
$ws.Range("L136").Value = 7016.400000000001
$ws.Range("M136").Value = -76925058
$ws.Range("N136").Value = -12116.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 2936.6667
$ws.Range("J132").Value = 3162.3076
$ws.Range("L132").Value = 28460.7684
$ws.Range("N132").Value = -33520.7684
# Row 137
$ws.Range("H137").Value = 38908.125
$ws.Range("I137").Value = 843.3333
$ws.Range("J137").Value = 47692.31
$ws.Range("K137").Value = 2529.9999
$ws.Range("L137").Value = 143076.93
$ws.Range("M137").Value = 2570.0001
$ws.Range("N137").Value = -153276.93

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 45840.4
$ws.Range("I132").Value = 30513.059
$ws.Range("J132").Value = 93215.82
$ws.Range("K132").Value = 91539.177
$ws.Range("L132").Value = 279647.46
$ws.Range("M132").Value = -89009.177
$ws.Range("N132").Value = -284707.46

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 36
$ws.Range("H36").Value = 57715
$ws.Range("J36").Value = 57715
$ws.Range("L36").Value = 57715
$ws.Range("N36").Value = -58839
# Row 122
$ws.Range("H122").Value = 3356.12
$ws.Range("I122").Value = 2872.0625
$ws.Range("J122").Value = 4216.6665
$ws.Range("K122").Value = 8616.1875
$ws.Range("L122").Value = 12649.9995
$ws.Range("M122").Value = -6166.1875
$ws.Range("N122").Value = -17549.9995
# Row 136
$ws.Range("H136").Value = 50345.414
$ws.Range("I136").Value = 35733.207
$ws.Range("J136").Value = 85658.25
$ws.Range("K136").Value = 107199.621
$ws.Range("L136").Value = 256974.75
$ws.Range("M136").Value = -104649.621
$ws.Range("N136").Value = -262074.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 13978.8
$ws.Range("J54").Value = 13978.8
$ws.Range("L54").Value = 13978.8
$ws.Range("N54").Value = -15018.8
# Row 81
$ws.Range("H81").Value = 32199.908
$ws.Range("I81").Value = 1432.2858
$ws.Range("J81").Value = 54870.79
$ws.Range("K81").Value = 2864.5716
$ws.Range("L81").Value = 109741.58
$ws.Range("M81").Value = -1803.5716
$ws.Range("N81").Value = -111863.58
# Row 84
$ws.Range("H84").Value = 32199.908
$ws.Range("I84").Value = 1432.2858
$ws.Range("J84").Value = 54870.79
$ws.Range("K84").Value = 14322.858
$ws.Range("L84").Value = 548707.9
$ws.Range("M84").Value = -9018.858
$ws.Range("N84").Value = -559315.9
# Row 126
$ws.Range("H126").Value = 1275
$ws.Range("I126").Value = 1536.8
$ws.Range("J126").Value = 838.6667
$ws.Range("K126").Value = 4610.4
$ws.Range("L126").Value = 2516.0001
$ws.Range("M126").Value = -2140.4
$ws.Range("N126").Value = -7456.0001
# Row 128
$ws.Range("H128").Value = 53000
$ws.Range("J128").Value = 53000
$ws.Range("L128").Value = 53000
$ws.Range("N128").Value = -62960
# Row 132
$ws.Range("H132").Value = 65088.97
$ws.Range("I132").Value = 56335.9
$ws.Range("J132").Value = 78555.234
$ws.Range("K132").Value = 169007.7
$ws.Range("L132").Value = 235665.702
$ws.Range("M132").Value = -166477.7
$ws.Range("N132").Value = -240725.702
# Row 136
$ws.Range("H136").Value = 40073.48
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 40073.48
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 120220.44
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -125320.44
